$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row below the existing data row (between row 2 and the
# rest of the - currently empty - sheet). Inserting here (rather than right
# above row 2) means the new row inherits formatting from row 2 above it
# (i.e. no special style), matching the target's unstyled data rows.
$ws.Rows("3:3").Insert()

# Duplicate the existing Honduras Liga Nacional row (row 2) down into the
# newly inserted row 3, preserving all of its original values/types/format.
$ws.Rows("2:2").Copy()
$ws.Rows("3:3").PasteSpecial(-4104)
$excel.CutCopyMode = $false

# Now overwrite row 2 with the new Algerian Ligue 1 match. Force the
# text-like columns to stay plain text (not auto-parsed as a date/time),
# then restore the "Normal" style so no stray formatting is left behind.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 5).NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "Algerian Ligue 1"
$ws.Cells.Item(2, 2).Value = "2025-12-25"
$ws.Cells.Item(2, 3).Value = "15:30:00"
$ws.Cells.Item(2, 4).Value = "Belouizdad"
$ws.Cells.Item(2, 5).Value = "ES Setif"
$ws.Range("A2:E2").Style = "Normal"

$ws.Cells.Item(2, 6).Value = 1.04
$ws.Cells.Item(2, 7).Value = 1000
$ws.Cells.Item(2, 8).Value = 1.04
$ws.Cells.Item(2, 9).Value = 1000
$ws.Cells.Item(2, 10).Value = 1.01
$ws.Cells.Item(2, 11).Value = 1000
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 1.24
$ws.Cells.Item(2, 17).Value = 1.01
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 0
$ws.Cells.Item(2, 22).Value = 0
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 0
$ws.Cells.Item(2, 25).Value = 0
$ws.Cells.Item(2, 26).Value = 0
$ws.Cells.Item(2, 27).Value = 0
$ws.Cells.Item(2, 28).Value = 0
$ws.Cells.Item(2, 29).Value = 0
$ws.Cells.Item(2, 30).Value = 0
$ws.Cells.Item(2, 31).Value = 0
$ws.Cells.Item(2, 32).Value = 0
$ws.Cells.Item(2, 33).Value = 0
$ws.Cells.Item(2, 34).Value = 0
$ws.Cells.Item(2, 35).Value = 0
$ws.Cells.Item(2, 36).Value = 0
$ws.Cells.Item(2, 37).Value = 0
$ws.Cells.Item(2, 38).Value = 0
$ws.Cells.Item(2, 39).Value = 0
$ws.Cells.Item(2, 40).Value = 0
$ws.Cells.Item(2, 41).Value = 0

# Row 3 (the duplicated Honduras Liga Nacional match) is otherwise unchanged
# from before the edit (K3 stays 950); only the newly written row 2 carries
# the updated 1000 value in the Odd_D_Lay column.
